# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1
    "G3"  = 0
    "G4"  = 1
    "G5"  = 1
    "G6"  = 2
    "G7"  = 1
    "G8"  = 2
    "G9"  = 2
    "G10" = 0
    "G11" = 2
    "G12" = 2
    "G13" = 2
    "G14" = 1
    "G15" = 2
    "G16" = 0
    "G17" = 0
    "G18" = 0
    "G19" = 2
    "G20" = 0
    "G22" = 1
    "G24" = 2
    "G25" = 1
    "G26" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
